# Update crypto price/volume data per Sat May 25 05:40:37 UTC 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.785.06"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.759.78"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.17"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.25"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.757.11"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("E10").Value = "  +4.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.33"
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.37"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.388.20"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.760.00"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.803.33"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.32"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.18"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "497.80"
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.73"
$ws.Range("E22").Value = "  +16.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.730"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.76"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000148"
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.37"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("E28").Value = "  +2.40%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.53"
$ws.Range("E30").Value = "  +5.16%  "
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("E32").Value = "  +3.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.13"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.907.22"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.694.17"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.02"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.85"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "444.68"
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.82"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.50"
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.66"
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.844.17"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.08"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("E51").Value = "  +2.86%  "
